$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# ===========================================================================
# 1. Structural edits (do these FIRST - whole-row inserts shift every
#    column, so nothing else should be staged in a scratch column before
#    this point).
#    a) insert a new row at 6            -> old rows 6-12 become 7-13
#    b) insert a new row at 11 (after the shift above, this lands right
#       after "Add spring constraint")   -> rows 11-13 become 12-14
# ===========================================================================
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(11).Insert()

# ===========================================================================
# 2. Now that no further row-shifting will happen, stash the formats that
#    are about to be overwritten by the status changes below, into a
#    scratch column (J) that nothing else will touch.
# ===========================================================================
$ws.Range("H7").Copy()  | Out-Null   # s7: TODO + box border
$ws.Range("J1").PasteSpecial($xlPasteFormats)
$ws.Range("H9").Copy()  | Out-Null   # s9: UNDERWAY, plain
$ws.Range("J2").PasteSpecial($xlPasteFormats)
$ws.Range("H10").Copy() | Out-Null   # s8: TODO, plain
$ws.Range("J3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ===========================================================================
# 3. New row 6 - "Change Performance timers" / "...broadphase and
#    narrowphase..." - status TODO (plain).
# ===========================================================================
$ws.Range("F6").Value = "Change Performance timers"
$ws.Range("G6").Value = "Add in broadphase and narrowphase like in the tutorials"
$ws.Range("H6").Value = "TODO"
$ws.Range("J3").Copy() | Out-Null
$ws.Range("H6").PasteSpecial($xlPasteFormats)

# ===========================================================================
# 4. Row 7 "TARGET PRACTISE" section header - status TODO -> UNDERWAY
#    (border variant).
# ===========================================================================
$ws.Range("H7").Value = "UNDERWAY"
$ws.Range("H3").Copy() | Out-Null
$ws.Range("H7").PasteSpecial($xlPasteFormats)

# ===========================================================================
# 5. Row 8 "Make Target Scene" - status TODO -> UNDERWAY (plain).
# ===========================================================================
$ws.Range("H8").Value = "UNDERWAY"
$ws.Range("J2").Copy() | Out-Null
$ws.Range("H8").PasteSpecial($xlPasteFormats)

# ===========================================================================
# 6. Row 9 "Use lambda to change colour" - status UNDERWAY -> DONE (new
#    green fill).
# ===========================================================================
$ws.Range("H9").Value = "DONE"
$ws.Range("H9").Interior.Color = 5287936   # RGB(0,176,80) = FF00B050

# ===========================================================================
# 7. Row 10 "Add spring constraint" - status TODO -> DONE.
# ===========================================================================
$ws.Range("H10").Value = "DONE"
$ws.Range("H10").Interior.Color = 5287936

# ===========================================================================
# 8. New row 11 - "Debug draw spring" / jagged-triangle note - status TODO
#    (plain).
# ===========================================================================
$ws.Range("F11").Value = "Debug draw spring"
$ws.Range("G11").Value = "Make the line several lines making up a jagged triagnle spring constraint"
$ws.Range("H11").Value = "TODO"
$ws.Range("J3").Copy() | Out-Null
$ws.Range("H11").PasteSpecial($xlPasteFormats)

# ===========================================================================
# 9. Rows 12-14 ("GPU ACCELERATION" / "Collision Response" / "Display
#    Number of Entities") keep their original text & styles untouched - the
#    row inserts above already shifted them into place.
# ===========================================================================

# ===========================================================================
# 10. New rows 15-16 - "SOFT BODY" section.
# ===========================================================================
$ws.Range("F15").Value = "SOFT BODY"
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F15").PasteSpecial($xlPasteFormats)
$ws.Range("G15").Value = ""
$ws.Range("G3").Copy() | Out-Null
$ws.Range("G15").PasteSpecial($xlPasteFormats)
$ws.Range("H15").Value = "TODO"
$ws.Range("J1").Copy() | Out-Null
$ws.Range("H15").PasteSpecial($xlPasteFormats)

$ws.Range("F16").Value = "Make spring constr generic"
$ws.Range("F13").Copy() | Out-Null
$ws.Range("F16").PasteSpecial($xlPasteFormats)
$ws.Range("G16").Value = "If necessary"
$ws.Range("H16").Value = "TODO"
$ws.Range("J3").Copy() | Out-Null
$ws.Range("H16").PasteSpecial($xlPasteFormats)

# ===========================================================================
# 11. Clean up scratch cells used to stash formats.
# ===========================================================================
$excel.CutCopyMode = $false
$ws.Range("J1:J3").Clear()

# ===========================================================================
# 12. Restore the selection shown in the saved file.
# ===========================================================================
$ws.Range("K20").Select()

Write-Host "Edit complete"
